$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A28").Value = "/Users/rayan/Downloads/d1/A008_B009_1025IF.RDC/d2 copy/efgh.R3D"
$ws.Range("B28").Value = "/Users/rayan/Downloads/d1_converted/A008_B009_1025IF.RDC/d2 copy/efgh.mov"
